# Fixed update to excel issue
#
# 1) Rename the "Requested quantity" headers on the first two sheets.
# 2) Add a new "PO Forecast" worksheet at the end of the workbook with
#    ds / PO_Forecast / yhat_lower / yhat_upper columns and 64 rows of
#    forecast data.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

# --- 1) Rename headers -----------------------------------------------
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the "PO Forecast" worksheet at the end --------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$headerValues = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($c = 1; $c -le 4; $c++) {
    $wsForecast.Cells.Item(1, $c).Value = $headerValues[$c - 1]
}
# Match the bold/centered/bordered header style used on the other sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$forecastData = @(
    @(44934.99999999999, 284, -348.4948783530393, 956.2349402613229),
    @(44941.99999999999, 291, -407.9809336932788, 950.9666698323473),
    @(44948.99999999999, 298, -345.659554512278, 1019.731574512816),
    @(44969.99999999999, 320, -313.3346092021697, 951.2861564750888),
    @(44976.99999999999, 327, -310.7250816535308, 966.3960968112427),
    @(44983.99999999999, 334, -382.0909257405734, 991.9548950719254),
    @(44990.99999999999, 341, -351.0081680501394, 957.0119042719376),
    @(45032.99999999999, 383, -251.0533597198169, 1025.607408626484),
    @(45046.99999999999, 398, -202.6411671693316, 1056.809400318921),
    @(45053.99999999999, 405, -339.4440953083469, 1124.334238844524),
    @(45060.99999999999, 412, -231.6671924094385, 1082.97662739009),
    @(45067.99999999999, 419, -297.5761251167341, 1089.097544798897),
    @(45088.99999999999, 440, -232.5140386986994, 1125.605856501795),
    @(45095.99999999999, 447, -216.306881405397, 1061.365670088623),
    @(45109.99999999999, 462, -173.6301284566047, 1116.041937349041),
    @(45116.99999999999, 469, -181.6099004250664, 1102.848531961411),
    @(45123.99999999999, 476, -213.3067834009479, 1145.851832563359),
    @(45144.99999999999, 497, -129.5554213377931, 1177.855206943175),
    @(45151.99999999999, 504, -176.0834243123007, 1156.373372564006),
    @(45158.99999999999, 511, -145.0456661316742, 1236.835797281022),
    @(45165.99999999999, 518, -117.0843148568582, 1114.194551460076),
    @(45172.99999999999, 525, -162.946893053724, 1164.47543835803),
    @(45186.99999999999, 540, -97.17033943672847, 1237.623872419402),
    @(45193.99999999999, 547, -140.5757021316517, 1212.029166248532),
    @(45207.99999999999, 561, -104.4008375870816, 1231.608228190406),
    @(45214.99999999999, 568, -105.4469393853772, 1260.050219481499),
    @(45221.99999999999, 575, -66.58688095112871, 1349.290785161501),
    @(45228.99999999999, 582, -82.1661224613915, 1220.448243374454),
    @(45235.99999999999, 589, -119.4515720362732, 1306.300237438945),
    @(45242.99999999999, 596, -38.53086904880346, 1281.399661772203),
    @(45256.99999999999, 611, -78.109585951604, 1244.088262661134),
    @(45263.99999999999, 618, -59.65807468633784, 1328.58783969121),
    @(45270.99999999999, 625, -50.70271349968225, 1255.511718689987),
    @(45277.99999999999, 632, -47.23501427102951, 1280.842109778268),
    @(45298.99999999999, 653, -24.56053054777377, 1310.536924944756),
    @(45305.99999999999, 660, -34.48939377884423, 1314.451007031424),
    @(45312.99999999999, 667, -37.73490832389096, 1286.136963288536),
    @(45319.99999999999, 674, 13.86180063624515, 1330.165250951165),
    @(45326.99999999999, 682, 39.62797581467976, 1332.059468580807),
    @(45333.99999999999, 689, 45.26245098478569, 1372.934517734351),
    @(45340.99999999999, 696, -4.232717506600593, 1348.373577487049),
    @(45347.99999999999, 703, 75.8533735726088, 1339.813172990242),
    @(45354.99999999999, 710, 29.83707425028582, 1395.969464366007),
    @(45361.99999999999, 717, 50.037752458412, 1360.895825199643),
    @(45396.99999999999, 753, 86.60206177155975, 1419.089896664846),
    @(45403.99999999999, 760, 101.371814637555, 1442.008268811863),
    @(45410.99999999999, 767, 105.6623024116737, 1419.966652230784),
    @(45417.99999999999, 774, 140.7669274490036, 1416.682439930681),
    @(45424.99999999999, 781, 140.5434194897718, 1464.968746267891),
    @(45445.99999999999, 802, 125.7449949049596, 1448.689508529765),
    @(45501.99999999999, 859, 149.2546686991991, 1468.042279408369),
    @(45508.99999999999, 866, 223.5019758313537, 1511.806182722167),
    @(45515.99999999999, 873, 239.4896677450776, 1558.168739583794),
    @(45564.99999999999, 923, 227.3137867840821, 1582.230975116474),
    @(45571.99999999999, 930, 266.2910074581663, 1607.065884819638),
    @(45599.99999999999, 958, 308.9416712056562, 1639.038032411338),
    @(45606.99999999999, 966, 277.4183286963977, 1594.194744383428),
    @(45613.99999999999, 973, 314.9721488625122, 1628.301883170902),
    @(45620.99999999999, 980, 269.6968642284992, 1690.662987633418),
    @(45627.99999999999, 987, 306.9037507649243, 1666.443373791295),
    @(45634.99999999999, 994, 295.6729423542998, 1676.06695568169),
    @(45641.99999999999, 1001, 332.7975101996922, 1662.273603557462),
    @(45648.99999999999, 1008, 380.7377752582626, 1671.375824646687),
    @(45655.99999999999, 1015, 319.1700507665872, 1699.065285595035)
)

$rowIndex = 2
foreach ($row in $forecastData) {
    $wsForecast.Cells.Item($rowIndex, 1).Value = $row[0]
    $wsForecast.Cells.Item($rowIndex, 2).Value = $row[1]
    $wsForecast.Cells.Item($rowIndex, 3).Value = $row[2]
    $wsForecast.Cells.Item($rowIndex, 4).Value = $row[3]

    # Match the date-formatted style used for column A on the other sheets
    $wsWeekly.Range("A2").Copy()
    $wsForecast.Range("A" + $rowIndex).PasteSpecial(-4122)

    $rowIndex++
}

$wb.Application.CutCopyMode = $false

# Restore the originally active sheet/tab so we don't introduce an
# unrelated change to the workbook's active-tab selection.
$wsWeekly.Activate()
